$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing "Related Requirement" value and fill in new columns ---
$ws.Range("B2").Value = "Car_SRS_28"
$ws.Range("C2").Value = "Functional"
$ws.Range("D2").Value = 'Validate "See more" button functionality'
$ws.Range("E2").Value = "1) open URL `"http://CarPurchasing`"`n2)Login with an existing account"
$ws.Range("F2").Value = "User name: customer`npassword: soso.soso123"
$ws.Range("G2").Value = '1)From home bage click on "see more" button at any car'
$ws.Range("H2").Value = "User should be redirct to car detaiLs  page witch contains information about the car"
$ws.Range("K2").Value = "passed"

# --- Row 3: move previous B2 value here, and fill in the rest of the new row ---
$ws.Range("B3").Value = "Car_SRS_20"
$ws.Range("C3").Value = "Functional"
$ws.Range("D3").Value = "Validate information in car details page"
$ws.Range("E3").Value = "1) open URL `"http://CarPurchasing`"`n2)Login with an existing account"
$ws.Range("F3").Value = "User name: customer`npassword: soso.soso123"
$ws.Range("G3").Value = '1)From home bage click on "see more" button at any car'
$ws.Range("H3").Value = "User should be redirct to car detaiLs page`nand this page contains ( the model year, Status, Avaliable Color, Price, and another information about it)"

# --- Style: change alignment of the "Related Requirement" style (used by B2) to left/top,
#     then copy that format onto B3 so both share the same style index ---
$ws.Range("B2").HorizontalAlignment = -4131 # xlLeft
$ws.Range("B2").VerticalAlignment = -4160   # xlTop
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row heights: header row shrinks a bit, all data rows become a uniform 73.5 ---
$ws.Rows.Item(1).RowHeight = 30.75
for ($r = 2; $r -le 20; $r++) {
    $ws.Rows.Item($r).RowHeight = 73.5
}

# --- sheet view: scroll so column F is the left-most visible column, and select F3 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("F3").Select()
